$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.114020294155487
$ws.Cells.Item(2, 4).Value = 1.106726162453675
$ws.Cells.Item(2, 5).Value = 1.124689888729501
$ws.Cells.Item(2, 6).Value = 1.12474854141054
$ws.Cells.Item(2, 9).Value = 1.062385059478128
$ws.Cells.Item(2, 10).Value = 1.118738945220518
$ws.Cells.Item(2, 11).Value = 1.109332640301551
$ws.Cells.Item(2, 12).Value = 1.127252547944958
$ws.Cells.Item(2, 13).Value = 1.127311060018601
$ws.Cells.Item(2, 14).Value = 1.120327682843882
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.116154103303256
$ws.Cells.Item(3, 4).Value = 1.108443021662417
$ws.Cells.Item(3, 5).Value = 1.126736621181659
$ws.Cells.Item(3, 6).Value = 1.126682965925106
$ws.Cells.Item(3, 9).Value = 1.062988382929346
$ws.Cells.Item(3, 10).Value = 1.120541114706407
$ws.Cells.Item(3, 11).Value = 1.110870494511624
$ws.Cells.Item(3, 12).Value = 1.12912274908591
$ws.Cells.Item(3, 13).Value = 1.129069212994311
$ws.Cells.Item(3, 14).Value = 1.122132411617153
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.117530364964778
$ws.Cells.Item(4, 4).Value = 1.109549621053423
$ws.Cells.Item(4, 5).Value = 1.128056818870262
$ws.Cells.Item(4, 6).Value = 1.127930436255279
$ws.Cells.Item(4, 9).Value = 1.063375207062208
$ws.Cells.Item(4, 10).Value = 1.121702505934145
$ws.Cells.Item(4, 11).Value = 1.111860789681805
$ws.Cells.Item(4, 12).Value = 1.130328297139676
$ws.Cells.Item(4, 13).Value = 1.130202180852437
$ws.Cells.Item(4, 14).Value = 1.123295452153648
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.118107904488581
$ws.Cells.Item(5, 4).Value = 1.11001382061201
$ws.Cells.Item(5, 5).Value = 1.128610854312303
$ws.Cells.Item(5, 6).Value = 1.128453882216181
$ws.Cells.Item(5, 9).Value = 1.063536982528922
$ws.Cells.Item(5, 10).Value = 1.12218964379811
$ws.Cells.Item(5, 11).Value = 1.11227597994416
$ws.Cells.Item(5, 12).Value = 1.130834031500845
$ws.Cells.Item(5, 13).Value = 1.130677382713492
$ws.Cells.Item(5, 14).Value = 1.123783281809255
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.118204815624428
$ws.Cells.Item(6, 4).Value = 1.110091702834643
$ws.Cells.Item(6, 5).Value = 1.128703822777554
$ws.Cells.Item(6, 6).Value = 1.128541713660249
$ws.Cells.Item(6, 9).Value = 1.063564096010718
$ws.Cells.Item(6, 10).Value = 1.122271371877358
$ws.Cells.Item(6, 11).Value = 1.112345626516352
$ws.Cells.Item(6, 12).Value = 1.130918883960291
$ws.Cells.Item(6, 13).Value = 1.130757107436175
$ws.Cells.Item(6, 14).Value = 1.12386512595175
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.117538086131252
$ws.Cells.Item(7, 4).Value = 1.109555827670876
$ws.Cells.Item(7, 5).Value = 1.128064225720385
$ws.Cells.Item(7, 6).Value = 1.12793743442761
$ws.Cells.Item(7, 9).Value = 1.063377372024307
$ws.Cells.Item(7, 10).Value = 1.121709019429774
$ws.Cells.Item(7, 11).Value = 1.111866341887939
$ws.Cells.Item(7, 12).Value = 1.130335058992449
$ws.Cells.Item(7, 13).Value = 1.130208534807837
$ws.Cells.Item(7, 14).Value = 1.123301974899187
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.114742363103522
$ws.Cells.Item(8, 4).Value = 1.107307290452678
$ws.Cells.Item(8, 5).Value = 1.125382470457564
$ws.Cells.Item(8, 6).Value = 1.125403179052992
$ws.Cells.Item(8, 9).Value = 1.062589699028146
$ws.Cells.Item(8, 10).Value = 1.119348990178656
$ws.Cells.Item(8, 11).Value = 1.109853370655136
$ws.Cells.Item(8, 12).Value = 1.127885556678794
$ws.Cells.Item(8, 13).Value = 1.127906216866098
$ws.Cells.Item(8, 14).Value = 1.120938594135833
$ws.Cells.Item(9, 2).Value = 1.019999999999999
$ws.Cells.Item(9, 3).Value = 1.109780652478722
$ws.Cells.Item(9, 4).Value = 1.103311056485884
$ws.Cells.Item(9, 5).Value = 1.120623820907326
$ws.Cells.Item(9, 6).Value = 1.120904070301054
$ws.Cells.Item(9, 9).Value = 1.061174022362222
$ws.Cells.Item(9, 10).Value = 1.115153070958129
$ws.Cells.Item(9, 11).Value = 1.106268635368665
$ws.Cells.Item(9, 12).Value = 1.123532993875453
$ws.Cells.Item(9, 13).Value = 1.123812473027094
$ws.Cells.Item(9, 14).Value = 1.116736716228734
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.106447479177156
$ws.Cells.Item(10, 4).Value = 1.100622739673747
$ws.Cells.Item(10, 5).Value = 1.117427646930096
$ws.Cells.Item(10, 6).Value = 1.11788075806661
$ws.Cells.Item(10, 9).Value = 1.060211079122534
$ws.Cells.Item(10, 10).Value = 1.112329345206791
$ws.Cells.Item(10, 11).Value = 1.10385229757691
$ws.Cells.Item(10, 12).Value = 1.120605495427844
$ws.Cells.Item(10, 13).Value = 1.121057235794144
$ws.Cells.Item(10, 14).Value = 1.113908980462943
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.104997794074492
$ws.Cells.Item(11, 4).Value = 1.099452649718866
$ws.Cells.Item(11, 5).Value = 1.116037696881896
$ws.Cells.Item(11, 6).Value = 1.116565638863721
$ws.Cells.Item(11, 9).Value = 1.059789448981736
$ws.Cells.Item(11, 10).Value = 1.111100052982923
$ws.Cells.Item(11, 11).Value = 1.102799434145458
$ws.Cells.Item(11, 12).Value = 1.119331417663715
$ws.Cells.Item(11, 13).Value = 1.119857698808902
$ws.Cells.Item(11, 14).Value = 1.112677942503116
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.104458323226205
$ws.Cells.Item(12, 4).Value = 1.099017095085261
$ws.Cells.Item(12, 5).Value = 1.115520478765699
$ws.Cells.Item(12, 6).Value = 1.116076214977354
$ws.Cells.Item(12, 9).Value = 1.059632124727107
$ws.Cells.Item(12, 10).Value = 1.110642420159152
$ws.Cells.Item(12, 11).Value = 1.102407342479804
$ws.Cells.Item(12, 12).Value = 1.118857171279976
$ws.Cells.Item(12, 13).Value = 1.119411134250378
$ws.Cells.Item(12, 14).Value = 1.112219659788244
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.104574087017813
$ws.Cells.Item(13, 4).Value = 1.099110565573896
$ws.Cells.Item(13, 5).Value = 1.115631466314977
$ws.Cells.Item(13, 6).Value = 1.116181240601221
$ws.Cells.Item(13, 9).Value = 1.05966590370611
$ws.Cells.Item(13, 10).Value = 1.110740630506346
$ws.Cells.Item(13, 11).Value = 1.102491493622026
$ws.Cells.Item(13, 12).Value = 1.118958944315495
$ws.Cells.Item(13, 13).Value = 1.119506969684988
$ws.Cells.Item(13, 14).Value = 1.112318009605396
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.104953221693999
$ws.Cells.Item(14, 4).Value = 1.099416665777815
$ws.Cells.Item(14, 5).Value = 1.115994962590181
$ws.Cells.Item(14, 6).Value = 1.116525202068207
$ws.Cells.Item(14, 9).Value = 1.059776459096798
$ws.Cells.Item(14, 10).Value = 1.111062245856945
$ws.Cells.Item(14, 11).Value = 1.102767044495152
$ws.Cells.Item(14, 12).Value = 1.119292236814573
$ws.Cells.Item(14, 13).Value = 1.119820806260511
$ws.Cells.Item(14, 14).Value = 1.112640081686683
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.10518668650385
$ws.Cells.Item(15, 4).Value = 1.099605140077553
$ws.Cells.Item(15, 5).Value = 1.116218800679506
$ws.Cells.Item(15, 6).Value = 1.116737003997614
$ws.Cells.Item(15, 9).Value = 1.05984448125092
$ws.Cells.Item(15, 10).Value = 1.111260267842723
$ws.Cells.Item(15, 11).Value = 1.102936685757448
$ws.Cells.Item(15, 12).Value = 1.119497456346596
$ws.Cells.Item(15, 13).Value = 1.120014037613942
$ws.Cells.Item(15, 14).Value = 1.112838384886395
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.106543552472651
$ws.Cells.Item(16, 4).Value = 1.100700265495343
$ws.Cells.Item(16, 5).Value = 1.117519764718027
$ws.Cells.Item(16, 6).Value = 1.117967909307976
$ws.Cells.Item(16, 9).Value = 1.060238962053214
$ws.Cells.Item(16, 10).Value = 1.112410787897837
$ws.Cells.Item(16, 11).Value = 1.103922032145387
$ws.Cells.Item(16, 12).Value = 1.120689913473477
$ws.Cells.Item(16, 13).Value = 1.121136705898534
$ws.Cells.Item(16, 14).Value = 1.113990538811952
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.107392943084684
$ws.Cells.Item(17, 4).Value = 1.101385576502549
$ws.Cells.Item(17, 5).Value = 1.118334201870555
$ws.Cells.Item(17, 6).Value = 1.118738396420638
$ws.Cells.Item(17, 9).Value = 1.0604851519343
$ws.Cells.Item(17, 10).Value = 1.113130692465072
$ws.Cells.Item(17, 11).Value = 1.104538337403109
$ws.Cells.Item(17, 12).Value = 1.121436163482057
$ws.Cells.Item(17, 13).Value = 1.121839167014228
$ws.Cells.Item(17, 14).Value = 1.114711465726247
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.107887763016189
$ws.Cells.Item(18, 4).Value = 1.10178472662019
$ws.Cells.Item(18, 5).Value = 1.118808673770934
$ws.Cells.Item(18, 6).Value = 1.119187231019692
$ws.Cells.Item(18, 9).Value = 1.060628300519983
$ws.Cells.Item(18, 10).Value = 1.11354996564819
$ws.Cells.Item(18, 11).Value = 1.104897185234324
$ws.Cells.Item(18, 12).Value = 1.121870817903115
$ws.Cells.Item(18, 13).Value = 1.122248275240576
$ws.Cells.Item(18, 14).Value = 1.115131334325377
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.108056380427283
$ws.Cells.Item(19, 4).Value = 1.101920728835596
$ws.Cells.Item(19, 5).Value = 1.118970359697329
$ws.Cells.Item(19, 6).Value = 1.11934017490255
$ws.Cells.Item(19, 9).Value = 1.060677034509566
$ws.Cells.Item(19, 10).Value = 1.113692820087627
$ws.Cells.Item(19, 11).Value = 1.105019436485257
$ws.Cells.Item(19, 12).Value = 1.122018919404444
$ws.Cells.Item(19, 13).Value = 1.122387665278045
$ws.Cells.Item(19, 14).Value = 1.115274391634503
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.107301875377339
$ws.Cells.Item(20, 4).Value = 1.101312109310054
$ws.Cells.Item(20, 5).Value = 1.118246880212153
$ws.Cells.Item(20, 6).Value = 1.118655790379252
$ws.Cells.Item(20, 9).Value = 1.060458784691408
$ws.Cells.Item(20, 10).Value = 1.113053519392931
$ws.Cells.Item(20, 11).Value = 1.104472279285284
$ws.Cells.Item(20, 12).Value = 1.121356162309235
$ws.Cells.Item(20, 13).Value = 1.121763864468631
$ws.Cells.Item(20, 14).Value = 1.114634183059492
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.104841603658445
$ws.Cells.Item(21, 4).Value = 1.099326552803604
$ws.Cells.Item(21, 5).Value = 1.115887947875061
$ws.Cells.Item(21, 6).Value = 1.116423939889197
$ws.Cells.Item(21, 9).Value = 1.059743923031085
$ws.Cells.Item(21, 10).Value = 1.110967566501952
$ws.Cells.Item(21, 11).Value = 1.102685929740531
$ws.Cells.Item(21, 12).Value = 1.119194118233596
$ws.Cells.Item(21, 13).Value = 1.11972841711595
$ws.Cells.Item(21, 14).Value = 1.112545267876145
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.103288970394421
$ws.Cells.Item(22, 4).Value = 1.098072754743831
$ws.Cells.Item(22, 5).Value = 1.114399404917147
$ws.Cells.Item(22, 6).Value = 1.115015291768355
$ws.Cells.Item(22, 9).Value = 1.059290336116143
$ws.Cells.Item(22, 10).Value = 1.109650135921621
$ws.Cells.Item(22, 11).Value = 1.101556917738563
$ws.Cells.Item(22, 12).Value = 1.11782897185642
$ws.Cells.Item(22, 13).Value = 1.118442833516603
$ws.Cells.Item(22, 14).Value = 1.111225966393279
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.104112607671182
$ws.Cells.Item(23, 4).Value = 1.098737937117706
$ws.Cells.Item(23, 5).Value = 1.1151890304426
$ws.Cells.Item(23, 6).Value = 1.115762563597196
$ws.Cells.Item(23, 9).Value = 1.059531185682045
$ws.Cells.Item(23, 10).Value = 1.110349100153506
$ws.Cells.Item(23, 11).Value = 1.102155992143919
$ws.Cells.Item(23, 12).Value = 1.118553219404306
$ws.Cells.Item(23, 13).Value = 1.119124906012298
$ws.Cells.Item(23, 14).Value = 1.111925923234546
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.10734302685438
$ws.Cells.Item(24, 4).Value = 1.101345307764032
$ws.Cells.Item(24, 5).Value = 1.118286338888276
$ws.Cells.Item(24, 6).Value = 1.118693118277228
$ws.Cells.Item(24, 9).Value = 1.060470700303083
$ws.Cells.Item(24, 10).Value = 1.113088392543539
$ws.Cells.Item(24, 11).Value = 1.104502130059267
$ws.Cells.Item(24, 12).Value = 1.12139231331234
$ws.Cells.Item(24, 13).Value = 1.121797892379458
$ws.Cells.Item(24, 14).Value = 1.114669105733974
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.111067723004478
$ws.Cells.Item(25, 4).Value = 1.104348341971569
$ws.Cells.Item(25, 5).Value = 1.121858117881187
$ws.Cells.Item(25, 6).Value = 1.122071305234231
$ws.Cells.Item(25, 9).Value = 1.061543345709455
$ws.Cells.Item(25, 10).Value = 1.116242373865494
$ws.Cells.Item(25, 11).Value = 1.107199958289348
$ws.Cells.Item(25, 12).Value = 1.124662675066351
$ws.Cells.Item(25, 13).Value = 1.124875299240053
$ws.Cells.Item(25, 14).Value = 1.117827566071172
